# Update leve profit-calculation cells across multiple crafting-class sheets
# (values refreshed from market data by scheduled runner)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2885.25
$ws.Range("I40").Value = 2049
$ws.Range("J40").Value = 3721.5
$ws.Range("K40").Value = 2049
$ws.Range("L40").Value = 3721.5
$ws.Range("M40").Value = -1874
$ws.Range("N40").Value = -4071.5
$ws.Range("H80").Value = 2533.3333
$ws.Range("J80").Value = 2600
$ws.Range("L80").Value = 7800
$ws.Range("N80").Value = -9796
$ws.Range("H83").Value = 2533.3333
$ws.Range("J83").Value = 2600
$ws.Range("L83").Value = 23400
$ws.Range("N83").Value = -33384
$ws.Range("H137").Value = 3118.6667
$ws.Range("I137").Value = 2419.7646
$ws.Range("K137").Value = 7259.293799999999
$ws.Range("M137").Value = -4709.293799999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 171
$ws.Range("I5").Value = 167.5
$ws.Range("K5").Value = 167.5
$ws.Range("M5").Value = -55.5
$ws.Range("H45").Value = 3874.75
$ws.Range("I45").Value = 3766.3333
$ws.Range("K45").Value = 3766.3333
$ws.Range("M45").Value = -3389.3333
$ws.Range("H63").Value = 990.5
$ws.Range("I63").Value = 978
$ws.Range("K63").Value = 978
$ws.Range("M63").Value = -292
$ws.Range("H66").Value = 990.5
$ws.Range("I66").Value = 978
$ws.Range("K66").Value = 4890
$ws.Range("M66").Value = -1458
$ws.Range("H80").Value = 20000
$ws.Range("I80").Value = 20000
$ws.Range("K80").Value = 20000
$ws.Range("M80").Value = -19002
$ws.Range("H83").Value = 20000
$ws.Range("I83").Value = 20000
$ws.Range("K83").Value = 60000
$ws.Range("M83").Value = -55008
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H141").Value = 146666.67
$ws.Range("I141").Value = 70000
$ws.Range("J141").Value = 185000
$ws.Range("K141").Value = 70000
$ws.Range("L141").Value = 185000
$ws.Range("M141").Value = -64820
$ws.Range("N141").Value = -195360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 171
$ws.Range("I4").Value = 167.5
$ws.Range("K4").Value = 167.5
$ws.Range("M4").Value = -52.5
$ws.Range("H22").Value = 224
$ws.Range("I22").Value = 224
$ws.Range("K22").Value = 224
$ws.Range("M22").Value = -51
$ws.Range("H50").Value = 60000
$ws.Range("J50").Value = 60000
$ws.Range("L50").Value = 60000
$ws.Range("N50").Value = -61148
$ws.Range("H80").Value = 1076.25
$ws.Range("I80").Value = 1052.8334
$ws.Range("J80").Value = 1146.5
$ws.Range("K80").Value = 1052.8334
$ws.Range("L80").Value = 1146.5
$ws.Range("M80").Value = -54.83339999999998
$ws.Range("N80").Value = -3142.5
$ws.Range("H83").Value = 1076.25
$ws.Range("I83").Value = 1052.8334
$ws.Range("J83").Value = 1146.5
$ws.Range("K83").Value = 5264.166999999999
$ws.Range("L83").Value = 5732.5
$ws.Range("M83").Value = -272.1669999999995
$ws.Range("N83").Value = -15716.5
$ws.Range("H86").Value = 5377.9165
$ws.Range("I86").Value = 5526.8
$ws.Range("J86").Value = 4633.5
$ws.Range("K86").Value = 5526.8
$ws.Range("L86").Value = 4633.5
$ws.Range("M86").Value = -4403.8
$ws.Range("N86").Value = -6879.5
$ws.Range("H89").Value = 5377.9165
$ws.Range("I89").Value = 5526.8
$ws.Range("J89").Value = 4633.5
$ws.Range("K89").Value = 27634
$ws.Range("L89").Value = 23167.5
$ws.Range("M89").Value = -22018
$ws.Range("N89").Value = -34399.5
$ws.Range("H96").Value = 19450
$ws.Range("I96").Value = 19450
$ws.Range("K96").Value = 19450
$ws.Range("M96").Value = -16704
$ws.Range("H99").Value = 1040.5
$ws.Range("I99").Value = 1040.5
$ws.Range("K99").Value = 1040.5
$ws.Range("M99").Value = 457.5
$ws.Range("H105").Value = 29449
$ws.Range("I105").Value = 829
$ws.Range("K105").Value = 829
$ws.Range("M105").Value = 918
$ws.Range("H107").Value = 999.5
$ws.Range("I107").Value = 999.5
$ws.Range("K107").Value = 999.5
$ws.Range("M107").Value = 920.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 5476.615
$ws.Range("I99").Value = 4511.6665
$ws.Range("J99").Value = 6303.7144
$ws.Range("K99").Value = 4511.6665
$ws.Range("L99").Value = 6303.7144
$ws.Range("M99").Value = -3013.6665
$ws.Range("N99").Value = -9299.714400000001
$ws.Range("H107").Value = 1096.2858
$ws.Range("I107").Value = 1180
$ws.Range("K107").Value = 1180
$ws.Range("M107").Value = 740
$ws.Range("H126").Value = 5476.615
$ws.Range("I126").Value = 4511.6665
$ws.Range("J126").Value = 6303.7144
$ws.Range("K126").Value = 13534.9995
$ws.Range("L126").Value = 18911.1432
$ws.Range("M126").Value = -11064.9995
$ws.Range("N126").Value = -23851.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 37.53846
$ws.Range("I2").Value = 30
$ws.Range("J2").Value = 49.6
$ws.Range("K2").Value = 180
$ws.Range("L2").Value = 297.6
$ws.Range("M2").Value = -67
$ws.Range("N2").Value = -523.6
$ws.Range("H34").Value = 1076.5385
$ws.Range("I34").Value = 550
$ws.Range("J34").Value = 1310.5555
$ws.Range("K34").Value = 1650
$ws.Range("L34").Value = 3931.6665
$ws.Range("M34").Value = -1566
$ws.Range("N34").Value = -4099.666499999999
$ws.Range("H46").Value = 979.4
$ws.Range("J46").Value = 979.4
$ws.Range("L46").Value = 2938.2
$ws.Range("N46").Value = -3120.2
$ws.Range("H63").Value = 33333
$ws.Range("J63").Value = 33333
$ws.Range("L63").Value = 99999
$ws.Range("N63").Value = -101497
$ws.Range("H66").Value = 33333
$ws.Range("J66").Value = 33333
$ws.Range("L66").Value = 299997
$ws.Range("N66").Value = -307485
$ws.Range("H81").Value = 3350
$ws.Range("I81").Value = 3350
$ws.Range("K81").Value = 10050
$ws.Range("M81").Value = -8927
$ws.Range("H84").Value = 3350
$ws.Range("I84").Value = 3350
$ws.Range("K84").Value = 30150
$ws.Range("M84").Value = -24534
$ws.Range("H88").Value = 15000
$ws.Range("J88").Value = 15000
$ws.Range("L88").Value = 45000
$ws.Range("N88").Value = -45856
$ws.Range("H91").Value = 15000
$ws.Range("J91").Value = 15000
$ws.Range("L91").Value = 45000
$ws.Range("N91").Value = -47964
$ws.Range("H103").Value = 747
$ws.Range("I103").Value = 653.3333
$ws.Range("J103").Value = 1028
$ws.Range("K103").Value = 1959.9999
$ws.Range("L103").Value = 3084
$ws.Range("M103").Value = -1080.9999
$ws.Range("N103").Value = -4842

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 83.55556
$ws.Range("I2").Value = 82.666664
$ws.Range("J2").Value = 84
$ws.Range("K2").Value = 82.666664
$ws.Range("L2").Value = 84
$ws.Range("M2").Value = 30.333336
$ws.Range("N2").Value = -310
$ws.Range("H4").Value = 10700.8
$ws.Range("I4").Value = 1750
$ws.Range("J4").Value = 16668
$ws.Range("K4").Value = 1750
$ws.Range("L4").Value = 16668
$ws.Range("M4").Value = -1638
$ws.Range("N4").Value = -16892
$ws.Range("H10").Value = 52111.223
$ws.Range("I10").Value = 89667.664
$ws.Range("J10").Value = 33333
$ws.Range("K10").Value = 89667.664
$ws.Range("L10").Value = 33333
$ws.Range("M10").Value = -89498.664
$ws.Range("N10").Value = -33671
$ws.Range("H58").Value = 4814.2
$ws.Range("I58").Value = 5720.5
$ws.Range("J58").Value = 4210
$ws.Range("K58").Value = 5720.5
$ws.Range("L58").Value = 4210
$ws.Range("M58").Value = -5443.5
$ws.Range("N58").Value = -4764
$ws.Range("H70").Value = 5572.25
$ws.Range("I70").Value = 5572.25
$ws.Range("K70").Value = 5572.25
$ws.Range("M70").Value = -5302.25
$ws.Range("H73").Value = 5572.25
$ws.Range("I73").Value = 5572.25
$ws.Range("K73").Value = 5572.25
$ws.Range("M73").Value = -4636.25
$ws.Range("H80").Value = 8158.4614
$ws.Range("I80").Value = 2276.25
$ws.Range("J80").Value = 17570
$ws.Range("K80").Value = 2276.25
$ws.Range("L80").Value = 17570
$ws.Range("M80").Value = -1278.25
$ws.Range("N80").Value = -19566
$ws.Range("H83").Value = 8158.4614
$ws.Range("I83").Value = 2276.25
$ws.Range("J83").Value = 17570
$ws.Range("K83").Value = 11381.25
$ws.Range("L83").Value = 87850
$ws.Range("M83").Value = -6389.25
$ws.Range("N83").Value = -97834

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5848
$ws.Range("I7").Value = 5031.3335
$ws.Range("J7").Value = 8298
$ws.Range("K7").Value = 5031.3335
$ws.Range("L7").Value = 8298
$ws.Range("M7").Value = -4919.3335
$ws.Range("N7").Value = -8522
$ws.Range("H82").Value = 1835.6
$ws.Range("I82").Value = 1822
$ws.Range("J82").Value = 1890
$ws.Range("K82").Value = 1822
$ws.Range("L82").Value = 1890
$ws.Range("M82").Value = -1461
$ws.Range("N82").Value = -2612
$ws.Range("H85").Value = 1835.6
$ws.Range("I85").Value = 1822
$ws.Range("J85").Value = 1890
$ws.Range("K85").Value = 1822
$ws.Range("L85").Value = 1890
$ws.Range("M85").Value = -574
$ws.Range("N85").Value = -4386
$ws.Range("H126").Value = 5848
$ws.Range("I126").Value = 5031.3335
$ws.Range("J126").Value = 8298
$ws.Range("K126").Value = 15094.0005
$ws.Range("L126").Value = 24894
$ws.Range("M126").Value = -12624.0005
$ws.Range("N126").Value = -29834

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 2481.5715
$ws.Range("I113").Value = 1517.75
$ws.Range("K113").Value = 4553.25
$ws.Range("M113").Value = -2383.25
$ws.Range("H132").Value = 1062.1765
$ws.Range("I132").Value = 1062.1765
$ws.Range("K132").Value = 3186.5295
$ws.Range("M132").Value = -656.5295000000001
